# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# price table with freshly scraped values, mirroring the automated
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Price values are prefixed with a leading apostrophe so Excel stores
# them as literal text (matching the original inline-string cells)
# instead of auto-parsing them as numbers, which would silently strip
# meaningful trailing zeros (e.g. "0.8720" -> 0.872) or reformat
# multi-dot thousand-grouped numbers (e.g. "27.595.89").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.595.89"
$ws.Range("D3").Value = "'1.836.49"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'314.41"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4284"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").Value = "'0.3662"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.07274"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'0.8720"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "'20.71"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'1.767.00"
$ws.Range("E12").Value = "  -8.40%  "
$ws.Range("D13").Value = "'5.426"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'0.06938"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'80.42"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'0.000008929"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'15.45"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'27.419.87"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").Value = "'5.161"
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("E23").Value = "  +4.75%  "
$ws.Range("D24").Value = "'2.013.85"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "'1.978"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'154.76"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'18.85"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'5.201"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "'114.68"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").Value = "'1.840"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'0.7600"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'4.549"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "'2.956"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'0.05313"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'0.5095"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'6.601"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'8.452"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "'10.51"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "'106.01"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'1.619"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "'1.765"
$ws.Range("E51").Value = "  +3.00%  "
